# Append " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document."), as three additional runs:
#   " (", "Changed main", ")"
#
# A plain sequence of Range.InsertAfter calls at the same insertion
# point gets coalesced into a single run on save (because the
# inserted text shares the same - empty - character formatting as the
# text it is appended to). To keep each piece of text in its own
# <w:r>, as the target document requires, each piece is first typed
# into its own throw-away paragraph (so it is unambiguously a
# separate run) and the intervening paragraph marks are then deleted
# to stitch everything back into paragraph 1. Joining paragraphs this
# way does not trigger run-coalescing, so the runs stay distinct.

$d = $word.ActiveDocument

$para1 = $d.Paragraphs(1).Range
$para1.MoveEnd(1, -1)          # exclude the trailing paragraph mark
$p1End = $para1.End            # insertion point right after "document."

$pieces = @(" (", "Changed main", ")")

$cursor = $p1End
$joinPositions = @()

foreach ($piece in $pieces) {
    $d.Range($cursor, $cursor).InsertParagraphAfter()
    $joinPositions += $cursor      # position of the paragraph mark to remove later
    $cursor = $cursor + 1          # step over the new paragraph mark
    $d.Range($cursor, $cursor).InsertAfter($piece)
    $cursor = $cursor + $piece.Length
}

# Remove the paragraph marks starting from the last one so that the
# positions recorded above stay valid for the earlier deletions.
for ($i = $joinPositions.Length - 1; $i -ge 0; $i--) {
    $pos = $joinPositions[$i]
    $d.Range($pos, $pos + 1).Delete()
}
